$wb = $excel.ActiveWorkbook

$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A2").Value = "Version: " + $newVersion
$aboutSheet.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Chaochuan Coal Mine No. 1 Well, China, M1834, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$boundSheet = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 8; $row++) {
    $boundSheet.Cells.Item($row, 19).Value = $newVersion
}
